$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the 16th-May-refresh rows (regcntr_id, usr_id, machine_id, lang_code,
# is_active, cr_by, cr_dtimes, eff_dtimes) to the bottom of the table.
$newRows = @(
    @(10005, 110033, 10005),
    @(10005, 110034, 10005),
    @(10005, 110035, 10005)
)

$r = 34
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = "eng"
    $ws.Cells.Item($r, 5).Value = $true
    $ws.Cells.Item($r, 6).Value = "superadmin"
    $ws.Cells.Item($r, 7).Value = "now()"
    $ws.Cells.Item($r, 8).Value = "now()"
    $r++
}

# Reflect the post-edit selection left behind in the saved file (cursor moved
# to the first empty row, whole-row-to-sheet-bottom selection).
$ws.Range("A37:XFD1048576").Select()
